$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (F1) and "is_enabled_lbl" (G1) columns entirely.
# This shifts the former H1 (order_by) and I1 (rem) columns left into F and G.
$ws.Range("F1:G1").EntireColumn.Delete()
